function Get-ParagraphForPosition($doc, $position) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $candidate = $doc.Paragraphs.Item($i)
        if ($position -ge $candidate.Range.Start -and $position -lt $candidate.Range.End) {
            return $candidate
        }
    }
    return $null
}

$d = $word.ActiveDocument

# --- Replace the "m:if self.oclIsKindOf(ecore::EClass)" field (a real Word
#     field, fldChar begin/instrText*/fldChar end) with plain literal text
#     runs reading "{m:if self.oclIsKindOf(ecore::EClass)}". This matches the
#     new TokenIteratorFieldRewriterSplit parser, which expects M2Doc
#     tokens as plain braced text instead of Word fields. ---
$ifField = $d.Fields.Item(1)
$ifPara = Get-ParagraphForPosition $d $ifField.Code.Start
$ifField.Delete()

$pos = $ifPara.Range.Start
$run = $d.Range($pos, $pos)
$run.InsertAfter("{m:if ")
$run = $d.Range($run.End, $run.End)
$run.InsertAfter("self.oclIsKindOf(ecore::EC")
$run = $d.Range($run.End, $run.End)
$run.InsertAfter("l")
$run = $d.Range($run.End, $run.End)
$run.InsertAfter("ass)}")

# --- Replace the "m:endif" field with literal text "{m:endif}". ---
$endifField = $d.Fields.Item(1)
$endifPara = Get-ParagraphForPosition $d $endifField.Code.Start
$endifField.Delete()

$pos2 = $endifPara.Range.Start
$run2 = $d.Range($pos2, $pos2)
$run2.InsertAfter("{m:endif}")
